# Chiffres COVID-19 Valais - update to 02.06.2020 edition
# - revises the "Patients COVID-19 aux SI total" (col G) historical series for rows 12-93
# - appends 4 new days of data (rows 95-98), shifting the former last row (94) into
#   the regular body and promoting a new row 98 to be the sheet's final (bottom-border) row
# - refreshes the title string and the default selection/scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Title shared string
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Données COVID-19 Valais 02.06.2020"

# ---------------------------------------------------------------------------
# 2) Revised column G ("Patients COVID-19 aux SI total") values for rows 12-93.
#    Column H (total hospitalisations) is a shared formula (=G+E) and recalcs
#    automatically once G changes.
# ---------------------------------------------------------------------------
$gvals = @(
    15,16,19,20,27,29,33,35,40,45,57,59,67,75,81,92,
    105,104,110,118,122,131,129,129,121,119,121,122,121,
    108,106,98,96,94,95,91,88,83,77,78,72,69,69,67,
    65,67,63,59,61,62,54,51,52,53,51,51,52,47,45,
    42,37,35,33,34,34,33,30,28,25,25,25,24,22,21,
    21,21,21,21,21,20,20,19,18
)
for ($i = 0; $i -lt $gvals.Length; $i++) {
    $row = 12 + $i
    $ws.Cells.Item($row, 7).Value = $gvals[$i]
}

# ---------------------------------------------------------------------------
# 3) Grow the table from row 94 down to row 98.
#    First, snapshot the current (pre-edit) formatting of the special last row
#    (row 94) onto the new last row (row 98), then stamp the regular interior
#    row formatting (row 93) onto rows 94-97.
# ---------------------------------------------------------------------------
$ws.Range("A94:L94").Copy() | Out-Null
$ws.Range("A98:L98").PasteSpecial(-4122) | Out-Null

$ws.Range("A93:L93").Copy() | Out-Null
$ws.Range("A94:L94").PasteSpecial(-4122) | Out-Null
$ws.Range("A95:L95").PasteSpecial(-4122) | Out-Null
$ws.Range("A96:L96").PasteSpecial(-4122) | Out-Null
$ws.Range("A97:L97").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Populate data rows 94-98.
#    Columns: A date, B cumulative positive cases (formula), C new cases,
#    D new hospital admissions, E ICU, F total-SI related, G SI patients,
#    H total hospitalisations (formula), I cumulative deaths (formula),
#    J new deaths (formula), K new hospital deaths, L new extra-hosp deaths.
# ---------------------------------------------------------------------------
$dates = @(43980, 43981, 43982, 43983, 43984)
$cVals = @(0, 0, 2, 1, 0)
$dVals = @(0, 0, 2, 1, 0)
$eVals = @(3, 3, 5, 4, 4)
$fVals = @(3, 3, 4, 4, 4)
$gVals2 = @(17, 16, 16, 17, 17)
$kVals = @(0, 0, 0, 0, 0)
$lVals = @(0, 0, 0, 0, 0)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 94 + $i

    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
    $ws.Cells.Item($row, 6).Value = $fVals[$i]
    $ws.Cells.Item($row, 7).Value = $gVals2[$i]
    $ws.Cells.Item($row, 11).Value = $kVals[$i]
    $ws.Cells.Item($row, 12).Value = $lVals[$i]

    $prev = $row - 1
    $ws.Cells.Item($row, 2).Formula = "=B$prev+C$row"
    $ws.Cells.Item($row, 8).Formula = "=G$row+E$row"
    $ws.Cells.Item($row, 9).Formula = "=I$prev+J$row"
    $ws.Cells.Item($row, 10).Formula = "=K$row+L$row"
}

# ---------------------------------------------------------------------------
# 5) Sheet view: drop the old frozen scroll position / D-column selection,
#    select the title row instead, and land the viewport back at A1.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1:L1").Select()
